$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B12").Value = "V případě, že je více možných reakcí, se názvy uloží za sebe."
$ws.Range("B13").Value = "Názvy mp3 souborů jsou pojmenované čísleně ve formátu 00X."
